$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)
    $cell.NumberFormat = "@"
    $cell.Value = "2008-06-23"
}

$wb.Save()
